# Updated cryptos list on Mon Jun 12 10:11:23 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.973.32'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '1.749.83'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = "'234.30"
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").Value = "'0.5197"
$ws.Range("E7").Value = '  +2.82%  '
$ws.Range("D8").Value = "'0.2839"
$ws.Range("E8").Value = '  +7.71%  '
$ws.Range("D9").Value = "'39.71"
$ws.Range("E9").Value = '  -2.55%  '
$ws.Range("D10").Value = "'0.06141"
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("D11").Value = '1.752.45'
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("D12").Value = "'0.07031"
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").Value = "'15.51"
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").Value = "'0.6453"
$ws.Range("E14").Value = '  +6.23%  '
$ws.Range("D15").Value = "'4.529"
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = "'77.44"
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("D17").Value = "'0.9996"
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = "'0.9996"
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = '25.964.07'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = "'11.52"
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("D21").Value = "'0.000006626"
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("D22").Value = '1.970.80'
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("D23").Value = "'4.161"
$ws.Range("E23").Value = '  +2.92%  '
$ws.Range("D24").Value = "'8.603"
$ws.Range("E24").Value = '  +5.12%  '
$ws.Range("D25").Value = "'5.163"
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").Value = "'138.88"
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("D27").Value = "'1.502"
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("D28").Value = "'1.856"
$ws.Range("E28").Value = '  +4.28%  '
$ws.Range("D29").Value = "'15.12"
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").Value = "'103.30"
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("D31").Value = "'0.08312"
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").Value = "'3.654"
$ws.Range("E32").Value = '  -1.19%  '
$ws.Range("D33").Value = "'3.449"
$ws.Range("E33").Value = '  +1.90%  '
$ws.Range("D34").Value = "'0.04432"
$ws.Range("E34").Value = '  +1.64%  '
$ws.Range("E35").Value = '  -1.74%  '
$ws.Range("D36").Value = "'0.9875"
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").Value = "'0.6111"
$ws.Range("E37").Value = '  +1.92%  '
$ws.Range("D38").Value = "'2.685"
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").Value = "'0.01580"
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("D40").Value = "'1.950"
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").Value = "'0.9991"
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("D43").Value = "'0.3878"
$ws.Range("E43").Value = '  +2.16%  '
$ws.Range("D44").Value = "'0.7340"
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("D45").Value = "'5.016"
$ws.Range("E45").Value = '  +3.53%  '
$ws.Range("D46").Value = "'0.05470"
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("D47").Value = "'6.381"
$ws.Range("E47").Value = '  +8.03%  '
$ws.Range("D48").Value = "'0.1124"
$ws.Range("E48").Value = '  +3.93%  '
$ws.Range("D49").Value = "'52.85"
$ws.Range("E49").Value = '  +1.49%  '
$ws.Range("D50").Value = "'30.07"
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = "'0.3433"
$ws.Range("E51").Value = '  +0.71%  '
